# Update cryptocurrency price and volume data for Sun Dec  3 19:57:50 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '39.534.58'
$ws.Range('E2').Value = '  +0.42%  '
$ws.Range('D3').Value = '2.166.22'
$ws.Range('E3').Value = '  +0.44%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '227.13'
$ws.Range('E5').Value = '  -0.79%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.622'
$ws.Range('E6').Value = '  +0.14%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '62.90'
$ws.Range('E7').Value = '  -0.23%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E9').Value = '  -0.96%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0853'
$ws.Range('E10').Value = '  -0.91%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.103'
$ws.Range('E11').Value = '  +0.12%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '15.84'
$ws.Range('E12').Value = '  -1.66%  '
$ws.Range('D13').Value = '2.487.52'
$ws.Range('E13').Value = '  +0.56%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '21.68'
$ws.Range('E14').Value = '  -2.77%  '
$ws.Range('E15').Value = '  -1.63%  '
$ws.Range('E16').Value = '  -1.50%  '
$ws.Range('D17').Value = '2.160.93'
$ws.Range('E17').Value = '  +0.76%  '
$ws.Range('D18').Value = '39.566.55'
$ws.Range('E18').Value = '  +0.21%  '
$ws.Range('D19').Value = '0.0₃0924'
$ws.Range('E19').Value = '  +7.97%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '71.62'
$ws.Range('E20').Value = '  -0.78%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.99'
$ws.Range('E21').Value = '  -2.53%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '227.47'
$ws.Range('E22').Value = '  -0.40%  '
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('E24').Value = '  +0.28%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.30'
$ws.Range('E25').Value = '  -3.22%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '170.31'
$ws.Range('E26').Value = '  -0.83%  '
$ws.Range('E27').Value = '  -3.36%  '
$ws.Range('E28').Value = '  -0.96%  '
$ws.Range('E29').Value = '  +1.91%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '19.65'
$ws.Range('E30').Value = '  -0.02%  '
$ws.Range('E31').Value = '  +4.26%  '
$ws.Range('E32').Value = '  +0.28%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.49'
$ws.Range('E33').Value = '  -3.12%  '
$ws.Range('E34').Value = '  -2.81%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '6.95'
$ws.Range('E35').Value = '  -3.04%  '
$ws.Range('E36').Value = '  -1.42%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.79'
$ws.Range('E37').Value = '  +6.67%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.38'
$ws.Range('E38').Value = '  -0.87%  '
$ws.Range('E39').Value = '  +0.37%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '4.88'
$ws.Range('E40').Value = '  +17.06%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '102.23'
$ws.Range('E41').Value = '  -0.75%  '
$ws.Range('E42').Value = '  -1.55%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '17.65'
$ws.Range('E43').Value = '  -3.25%  '
$ws.Range('D44').Value = '1.513.14'
$ws.Range('E44').Value = '  -1.48%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.20'
$ws.Range('E45').Value = '  +0.92%  '
$ws.Range('E46').Value = '  +0.83%  '
$ws.Range('E47').Value = '  -0.36%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0914'
$ws.Range('E48').Value = '  -0.83%  '
$ws.Range('E49').Value = '  -1.60%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.000194'
$ws.Range('E50').Value = '  +29.78%  '
$ws.Range('D51').Value = '2.369.14'
$ws.Range('E51').Value = '  +0.48%  '
